$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 39 (pushes old rows 39..173 down to 40..174,
# and grows the used range to R174).
$ws.Rows.Item(39).Insert()

# Populate the freshly inserted row 39 with this week's new data point.
# (Same mercado/región/categoría/calidad/unidad/origen/clasificación as
# the row that used to be here; only the date + volume + price fields
# are new.)
$ws.Range("A39").Value = 7
$ws.Range("B39").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C39").Value = "Ñuble"
$ws.Range("D39").Value = 44487
$ws.Range("E39").Value = 16
$ws.Range("F39").Value = 100112023
$ws.Range("G39").Value = "Brócoli"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 750
$ws.Range("L39").Value = 800
$ws.Range("M39").Value = 775
$ws.Range("N39").Value = "$/unidad"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 775
$ws.Range("Q39").Value = 1
$ws.Range("R39").Value = "Hortaliza"
